$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2032520325203252
$ws.Cells.Item(2, 3).Value = 0.5447154471544715
$ws.Cells.Item(2, 10).Value = 0.03658536585365853
$ws.Cells.Item(2, 15).Value = 0.004065040650406504
$ws.Cells.Item(2, 16).Value = 0.1138211382113821
$ws.Cells.Item(2, 19).Value = 0.0975609756097561
$ws.Cells.Item(3, 2).Value = 0.01351351351351351
$ws.Cells.Item(3, 3).Value = 0.0472972972972973
$ws.Cells.Item(3, 10).Value = 0.04054054054054054
$ws.Cells.Item(3, 16).Value = 0.75
$ws.Cells.Item(3, 19).Value = 0.1486486486486487
$ws.Cells.Item(4, 10).Value = 0.06060606060606061
$ws.Cells.Item(4, 16).Value = 0.7272727272727273
$ws.Cells.Item(4, 19).Value = 0.2121212121212121
$ws.Cells.Item(5, 10).Value = 0.2
$ws.Cells.Item(5, 16).Value = 0.6
$ws.Cells.Item(5, 19).Value = 0.2
$ws.Cells.Item(6, 2).Value = 0.02325581395348837
$ws.Cells.Item(6, 5).Value = 0.004651162790697674
$ws.Cells.Item(6, 6).Value = 0.07441860465116279
$ws.Cells.Item(6, 10).Value = 0.2790697674418605
$ws.Cells.Item(6, 15).Value = 0.02325581395348837
$ws.Cells.Item(6, 17).Value = 0.1813953488372093
$ws.Cells.Item(6, 18).Value = 0.04186046511627907
$ws.Cells.Item(6, 19).Value = 0.3720930232558139
$ws.Cells.Item(7, 2).Value = 0.107981220657277
$ws.Cells.Item(7, 4).Value = 0.01408450704225352
$ws.Cells.Item(7, 5).Value = 0.004694835680751174
$ws.Cells.Item(7, 6).Value = 0.04225352112676056
$ws.Cells.Item(7, 10).Value = 0.1784037558685446
$ws.Cells.Item(7, 15).Value = 0.0187793427230047
$ws.Cells.Item(7, 17).Value = 0.1737089201877934
$ws.Cells.Item(7, 18).Value = 0.07981220657276995
$ws.Cells.Item(7, 19).Value = 0.3802816901408451
$ws.Cells.Item(8, 2).Value = 0.06060606060606061
$ws.Cells.Item(8, 4).Value = 0.01864801864801865
$ws.Cells.Item(8, 6).Value = 0.0675990675990676
$ws.Cells.Item(8, 10).Value = 0.1421911421911422
$ws.Cells.Item(8, 15).Value = 0.02097902097902098
$ws.Cells.Item(8, 17).Value = 0.1771561771561772
$ws.Cells.Item(8, 18).Value = 0.1282051282051282
$ws.Cells.Item(8, 19).Value = 0.3846153846153846
$ws.Cells.Item(9, 2).Value = 0.08860759493670886
$ws.Cells.Item(9, 4).Value = 0.01687763713080169
$ws.Cells.Item(9, 5).Value = 0.008438818565400843
$ws.Cells.Item(9, 6).Value = 0.02953586497890295
$ws.Cells.Item(9, 10).Value = 0.1223628691983122
$ws.Cells.Item(9, 15).Value = 0.01687763713080169
$ws.Cells.Item(9, 17).Value = 0.1645569620253164
$ws.Cells.Item(9, 18).Value = 0.08438818565400844
$ws.Cells.Item(9, 19).Value = 0.4683544303797468
$ws.Cells.Item(10, 2).Value = 0.0936
$ws.Cells.Item(10, 4).Value = 0.0144
$ws.Cells.Item(10, 5).Value = 0.0016
$ws.Cells.Item(10, 6).Value = 0.068
$ws.Cells.Item(10, 10).Value = 0.144
$ws.Cells.Item(10, 15).Value = 0.0104
$ws.Cells.Item(10, 17).Value = 0.2144
$ws.Cells.Item(10, 18).Value = 0.0776
$ws.Cells.Item(10, 19).Value = 0.376
$ws.Cells.Item(11, 7).Value = 0.1118881118881119
$ws.Cells.Item(11, 10).Value = 0.06643356643356643
$ws.Cells.Item(11, 11).Value = 0.1643356643356643
$ws.Cells.Item(11, 12).Value = 0.6258741258741258
$ws.Cells.Item(11, 19).Value = 0.03146853146853147
$ws.Cells.Item(12, 7).Value = 0.8115183246073299
$ws.Cells.Item(12, 10).Value = 0.1204188481675393
$ws.Cells.Item(12, 11).Value = 0.01047120418848168
$ws.Cells.Item(12, 12).Value = 0.03664921465968586
$ws.Cells.Item(12, 19).Value = 0.02094240837696335
$ws.Cells.Item(13, 7).Value = 0.7272727272727273
$ws.Cells.Item(13, 10).Value = 0.2272727272727273
$ws.Cells.Item(13, 19).Value = 0.04545454545454546
$ws.Cells.Item(15, 6).Value = 0.0425531914893617
$ws.Cells.Item(15, 8).Value = 0.1106382978723404
$ws.Cells.Item(15, 9).Value = 0.09787234042553192
$ws.Cells.Item(15, 10).Value = 0.3276595744680851
$ws.Cells.Item(15, 11).Value = 0.08936170212765958
$ws.Cells.Item(15, 13).Value = 0.01276595744680851
$ws.Cells.Item(15, 15).Value = 0.09361702127659574
$ws.Cells.Item(15, 19).Value = 0.225531914893617
$ws.Cells.Item(16, 6).Value = 0.006329113924050633
$ws.Cells.Item(16, 8).Value = 0.2088607594936709
$ws.Cells.Item(16, 9).Value = 0.1012658227848101
$ws.Cells.Item(16, 10).Value = 0.3734177215189873
$ws.Cells.Item(16, 11).Value = 0.0949367088607595
$ws.Cells.Item(16, 13).Value = 0.0379746835443038
$ws.Cells.Item(16, 14).Value = 0.006329113924050633
$ws.Cells.Item(16, 15).Value = 0.05696202531645569
$ws.Cells.Item(16, 19).Value = 0.1139240506329114
$ws.Cells.Item(17, 6).Value = 0.01986754966887417
$ws.Cells.Item(17, 8).Value = 0.1832229580573951
$ws.Cells.Item(17, 9).Value = 0.1169977924944812
$ws.Cells.Item(17, 10).Value = 0.3686534216335541
$ws.Cells.Item(17, 11).Value = 0.09713024282560706
$ws.Cells.Item(17, 13).Value = 0.02869757174392936
$ws.Cells.Item(17, 14).Value = 0.002207505518763797
$ws.Cells.Item(17, 15).Value = 0.0728476821192053
$ws.Cells.Item(17, 19).Value = 0.1103752759381898
$ws.Cells.Item(18, 6).Value = 0.02564102564102564
$ws.Cells.Item(18, 8).Value = 0.2102564102564103
$ws.Cells.Item(18, 9).Value = 0.1282051282051282
$ws.Cells.Item(18, 10).Value = 0.3794871794871795
$ws.Cells.Item(18, 11).Value = 0.08205128205128205
$ws.Cells.Item(18, 13).Value = 0.01538461538461539
$ws.Cells.Item(18, 15).Value = 0.04615384615384616
$ws.Cells.Item(18, 19).Value = 0.1128205128205128
$ws.Cells.Item(19, 6).Value = 0.01050080775444265
$ws.Cells.Item(19, 8).Value = 0.1954765751211632
$ws.Cells.Item(19, 9).Value = 0.09693053311793215
$ws.Cells.Item(19, 10).Value = 0.3796445880452343
$ws.Cells.Item(19, 11).Value = 0.1098546042003231
$ws.Cells.Item(19, 13).Value = 0.01696284329563813
$ws.Cells.Item(19, 14).Value = 0.0008077544426494346
$ws.Cells.Item(19, 15).Value = 0.0815831987075929
$ws.Cells.Item(19, 19).Value = 0.1082390953150242
